$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 "30.391.38"
Set-TextValue 2 5 "  +0.04%  "

# Row 3
Set-TextValue 3 4 "1.877.43"
Set-TextValue 3 5 "  -0.74%  "

# Row 5
Set-TextValue 5 4 "238.58"
Set-TextValue 5 5 "  +0.16%  "

# Row 6
Set-TextValue 6 5 "  +0.02%  "

# Row 7
Set-TextValue 7 4 "0.4797"
Set-TextValue 7 5 "  -0.70%  "

# Row 8
Set-TextValue 8 4 "0.2820"
Set-TextValue 8 5 "  -2.80%  "

# Row 9
Set-TextValue 9 4 "0.06512"
Set-TextValue 9 5 "  -1.42%  "

# Row 10
Set-TextValue 10 4 "1.873.84"
Set-TextValue 10 5 "  -0.95%  "

# Row 11
Set-TextValue 11 4 "0.07477"
Set-TextValue 11 5 "  +0.77%  "

# Row 12
Set-TextValue 12 4 "16.61"
Set-TextValue 12 5 "  -1.64%  "

# Row 13
Set-TextValue 13 4 "5.091"
Set-TextValue 13 5 "  -1.42%  "

# Row 14
Set-TextValue 14 4 "88.16"
Set-TextValue 14 5 "  +0.70%  "

# Row 15
Set-TextValue 15 4 "0.6605"
Set-TextValue 15 5 "  -0.38%  "

# Row 16
Set-TextValue 16 4 "30.364.94"
Set-TextValue 16 5 "  +0.06%  "

# Row 17
Set-TextValue 17 4 "13.30"
Set-TextValue 17 5 "  -0.84%  "

# Row 18
Set-TextValue 18 4 "1.000"
Set-TextValue 18 5 "  +0.01%  "

# Row 19
Set-TextValue 19 4 "0.000007591"
Set-TextValue 19 5 "  -2.31%  "

# Row 20
Set-TextValue 20 4 "2.116.72"
Set-TextValue 20 5 "  -1.02%  "

# Row 21
Set-TextValue 21 4 "5.299"
Set-TextValue 21 5 "  -2.08%  "

# Row 22
Set-TextValue 22 4 "1.001"
Set-TextValue 22 5 "  +0.08%  "

# Row 23
Set-TextValue 23 4 "220.50"
Set-TextValue 23 5 "  +14.04%  "

# Row 24
Set-TextValue 24 4 "6.191"
Set-TextValue 24 5 "  +0.30%  "

# Row 25
Set-TextValue 25 4 "9.361"
Set-TextValue 25 5 "  -0.56%  "

# Row 26
Set-TextValue 26 4 "168.08"
Set-TextValue 26 5 "  +2.92%  "

# Row 27
Set-TextValue 27 4 "18.41"
Set-TextValue 27 5 "  +0.88%  "

# Row 28
Set-TextValue 28 4 "1.969"
Set-TextValue 28 5 "  +1.02%  "

# Row 29
Set-TextValue 29 4 "1.463"
Set-TextValue 29 5 "  +1.16%  "

# Row 30
Set-TextValue 30 2 "Stellar"
Set-TextValue 30 3 "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue 30 4 "0.09356"
Set-TextValue 30 5 "  +2.49%  "

# Row 31
Set-TextValue 31 2 "InternetComputer(DFINITY)"
Set-TextValue 31 3 "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue 31 4 "4.314"
Set-TextValue 31 5 "  +0.34%  "

# Row 32
Set-TextValue 32 4 "4.027"
Set-TextValue 32 5 "  -0.49%  "

# Row 33
Set-TextValue 33 4 "0.05038"
Set-TextValue 33 5 "  -1.47%  "

# Row 34
Set-TextValue 34 4 "1.199"
Set-TextValue 34 5 "  +4.09%  "

# Row 35
Set-TextValue 35 4 "0.7433"
Set-TextValue 35 5 "  +1.36%  "

# Row 36
Set-TextValue 36 4 "2.710"
Set-TextValue 36 5 "  +0.06%  "

# Row 37
Set-TextValue 37 4 "0.01822"
Set-TextValue 37 5 "  +1.38%  "

# Row 38
Set-TextValue 38 5 "  -1.20%  "

# Row 39
Set-TextValue 39 4 "2.062"
Set-TextValue 39 5 "  -1.01%  "

# Row 40
Set-TextValue 40 4 "0.9044"
Set-TextValue 40 5 "  -1.47%  "

# Row 41
Set-TextValue 41 4 "106.59"
Set-TextValue 41 5 "  -0.06%  "

# Row 42
Set-TextValue 42 4 "5.883"
Set-TextValue 42 5 "  -0.52%  "

# Row 43
Set-TextValue 43 4 "0.4275"
Set-TextValue 43 5 "  -1.18%  "

# Row 44
Set-TextValue 44 4 "1.005"
Set-TextValue 44 5 "  +0.29%  "

# Row 45
Set-TextValue 45 4 "7.411"
Set-TextValue 45 5 "  -3.27%  "

# Row 46
Set-TextValue 46 4 "64.82"
Set-TextValue 46 5 "  -0.04%  "

# Row 47
Set-TextValue 47 4 "0.1276"
Set-TextValue 47 5 "  -4.18%  "

# Row 48
Set-TextValue 48 4 "1.481"
Set-TextValue 48 5 "  -5.54%  "

# Row 49
Set-TextValue 49 4 "8.937"
Set-TextValue 49 5 "  -0.35%  "

# Row 50
Set-TextValue 50 4 "33.74"
Set-TextValue 50 5 "  -0.78%  "

# Row 51
Set-TextValue 51 4 "0.3887"
Set-TextValue 51 5 "  +0.38%  "
